# Add working set of sequences
# Fill columns G:N (image1, type1, image2, type2, image3, type3, image4, type4)
# with "N/A" for the rows that only had data through column F.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(2,6,8,10,13,16,17,18,23,28,30,31,33,40,42,47,48,49,50,53,55,58,62,65,67,72,79,81,88,90,91,92,96,97,100,101,104,106,107,111,114,116,117,120,125,126,130,137,138,141,142,144,146,147,154,157,160,162,166,169,175,179,182,188)

foreach ($r in $rows) {
    $ws.Range("G" + $r + ":N" + $r).Value = "N/A"
}
